# Update routing/region files and add new license numbers to the ltcf
# comprehensive list.
#
# The "routes" sheet maps zip codes to a sub-region label. Rows 48, 49 and
# 75 had their region values corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("routes")

# Row 48: zip 98056 should be region "east_king_county" (was "south_king_county")
$ws.Range("C48").Value = "east_king_county"

# Row 49: zip 98057 should be region "south_king_county" (was "east_king_county")
$ws.Range("C49").Value = "south_king_county"

# Row 75: zip 98126 should be region "west_king_county" (was "south_king_county")
$ws.Range("C75").Value = "west_king_county"

# Update the sheet's view/scroll position and active selection to match
# where the editor was working when the change was saved.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 71
$ws.Range("E76").Select()
